$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.811.31'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').Value = '2.213.13'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '261.21'
$ws.Range('E5').Value = '  +2.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '86.49'
$ws.Range('E6').Value = '  +13.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.619'
$ws.Range('E7').Value = '  +1.49%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.601'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.15'
$ws.Range('E10').Value = '  +8.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0919'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.43'
$ws.Range('E12').Value = '  +7.34%  '
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('D14').Value = '2.547.22'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.46'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '2.211.78'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.784'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '43.748.44'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.95'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.82'
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.35'
$ws.Range('E22').Value = '  +7.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '231.58'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.90'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.55'
$ws.Range('E26').Value = '  +5.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.66'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.92'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('E29').Value = '  +3.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.68'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.46'
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0874'
$ws.Range('E33').Value = '  +3.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.42'
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('E37').Value = '  +4.67%  '
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.96'
$ws.Range('E39').Value = '  +8.01%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.60'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.53'
$ws.Range('E42').Value = '  +6.42%  '
$ws.Range('E43').Value = '  +4.74%  '
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.62'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0979'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.32'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.18'
$ws.Range('E48').Value = '  +4.36%  '
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.444'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.50'
$ws.Range('E51').Value = '  +5.45%  '
